# chore: update Sheets via scheduled runner
# Refreshes cached market-board figures (currentAveragePrice / NQ / HQ,
# LevePrice NQ/HQ, LeveProfit NQ/HQ -- columns H:N) across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW and WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2440
$ws.Range("I113").Value = 2620
$ws.Range("J113").Value = 1990
$ws.Range("K113").Value = 2620
$ws.Range("L113").Value = 1990
$ws.Range("M113").Value = 634
$ws.Range("N113").Value = -8498

$ws.Range("H129").Value = 889.4269399999999
$ws.Range("I129").Value = 451.08334
$ws.Range("J129").Value = 957.74023
$ws.Range("K129").Value = 1353.25002
$ws.Range("L129").Value = 2873.22069
$ws.Range("M129").Value = 3646.74998
$ws.Range("N129").Value = -12873.22069

$ws.Range("H132").Value = 1141.2424
$ws.Range("I132").Value = 888.5
$ws.Range("K132").Value = 2665.5
$ws.Range("M132").Value = -135.5

$ws.Range("H133").Value = 35156
$ws.Range("J133").Value = 35156
$ws.Range("L133").Value = 35156
$ws.Range("N133").Value = -45276

$ws.Range("H137").Value = 1579
$ws.Range("I137").Value = 1324
$ws.Range("J137").Value = 2344
$ws.Range("K137").Value = 3972
$ws.Range("L137").Value = 7032
$ws.Range("M137").Value = -1422
$ws.Range("N137").Value = -12132

$ws.Range("H138").Value = 5359.1353
$ws.Range("I138").Value = 883.4815
$ws.Range("J138").Value = 17443.4
$ws.Range("K138").Value = 2650.4445
$ws.Range("L138").Value = 52330.2
$ws.Range("M138").Value = 2489.5555
$ws.Range("N138").Value = -62610.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3652.818
$ws.Range("I2").Value = 3522.5
$ws.Range("J2").Value = 4000.3333
$ws.Range("K2").Value = 3522.5
$ws.Range("L2").Value = 4000.3333
$ws.Range("M2").Value = -3409.5
$ws.Range("N2").Value = -4226.3333

$ws.Range("H32").Value = 5046.1387
$ws.Range("I32").Value = 3948.5933
$ws.Range("J32").Value = 10027.308
$ws.Range("K32").Value = 3948.5933
$ws.Range("L32").Value = 10027.308
$ws.Range("M32").Value = -3661.5933
$ws.Range("N32").Value = -10601.308

$ws.Range("H45").Value = 5568.522
$ws.Range("I45").Value = 9483
$ws.Range("K45").Value = 9483
$ws.Range("M45").Value = -9106

$ws.Range("H61").Value = 5233.759
$ws.Range("I61").Value = 5477
$ws.Range("J61").Value = 1950
$ws.Range("K61").Value = 5477
$ws.Range("L61").Value = 1950
$ws.Range("M61").Value = -5265
$ws.Range("N61").Value = -2374

$ws.Range("H116").Value = 3652.818
$ws.Range("I116").Value = 3522.5
$ws.Range("J116").Value = 4000.3333
$ws.Range("K116").Value = 3522.5
$ws.Range("L116").Value = 4000.3333
$ws.Range("M116").Value = -1228.5
$ws.Range("N116").Value = -8588.3333

$ws.Range("H132").Value = 4455.92
$ws.Range("I132").Value = 2045.6666
$ws.Range("J132").Value = 6680.769
$ws.Range("K132").Value = 6136.9998
$ws.Range("L132").Value = 20042.307
$ws.Range("M132").Value = -3606.9998
$ws.Range("N132").Value = -25102.307

$ws.Range("H136").Value = 5233.759
$ws.Range("I136").Value = 5477
$ws.Range("J136").Value = 1950
$ws.Range("K136").Value = 16431
$ws.Range("L136").Value = 5850
$ws.Range("M136").Value = -13881
$ws.Range("N136").Value = -10950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3652.818
$ws.Range("I3").Value = 3522.5
$ws.Range("J3").Value = 4000.3333
$ws.Range("K3").Value = 3522.5
$ws.Range("L3").Value = 4000.3333
$ws.Range("M3").Value = -3408.5
$ws.Range("N3").Value = -4228.3333

$ws.Range("H94").Value = 1791.0555
$ws.Range("I94").Value = 1167.1818
$ws.Range("K94").Value = 1167.1818
$ws.Range("M94").Value = -716.1818000000001

$ws.Range("H122").Value = 20999
$ws.Range("J122").Value = 20999
$ws.Range("L122").Value = 20999
$ws.Range("N122").Value = -30799

$ws.Range("H133").Value = 31165.455
$ws.Range("J133").Value = 31165.455
$ws.Range("L133").Value = 31165.455
$ws.Range("N133").Value = -41285.455

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6904.36
$ws.Range("I31").Value = 1646.5294
$ws.Range("J31").Value = 18077.25
$ws.Range("K31").Value = 1646.5294
$ws.Range("L31").Value = 18077.25
$ws.Range("M31").Value = -1351.5294
$ws.Range("N31").Value = -18667.25

$ws.Range("H34").Value = 6904.36
$ws.Range("I34").Value = 1646.5294
$ws.Range("J34").Value = 18077.25
$ws.Range("K34").Value = 1646.5294
$ws.Range("L34").Value = 18077.25
$ws.Range("M34").Value = -1444.5294
$ws.Range("N34").Value = -18481.25

$ws.Range("H58").Value = 1661.9231
$ws.Range("I58").Value = 925.41174
$ws.Range("J58").Value = 2231.0454
$ws.Range("K58").Value = 925.41174
$ws.Range("L58").Value = 2231.0454
$ws.Range("M58").Value = -722.41174
$ws.Range("N58").Value = -2637.0454

$ws.Range("H81").Value = 33147.367
$ws.Range("J81").Value = 33147.367
$ws.Range("L81").Value = 33147.367
$ws.Range("N81").Value = -35143.367

$ws.Range("H84").Value = 33147.367
$ws.Range("J84").Value = 33147.367
$ws.Range("L84").Value = 99442.101
$ws.Range("N84").Value = -109426.101

$ws.Range("H94").Value = 2816.1333
$ws.Range("J94").Value = 2390.9546
$ws.Range("L94").Value = 2390.9546
$ws.Range("N94").Value = -3292.9546

$ws.Range("H99").Value = 5437720.5
$ws.Range("I99").Value = 1792.1333
$ws.Range("J99").Value = 15630086
$ws.Range("K99").Value = 1792.1333
$ws.Range("L99").Value = 15630086
$ws.Range("M99").Value = -294.1333
$ws.Range("N99").Value = -15633082

$ws.Range("H126").Value = 5437720.5
$ws.Range("I126").Value = 1792.1333
$ws.Range("J126").Value = 15630086
$ws.Range("K126").Value = 5376.3999
$ws.Range("L126").Value = 46890258
$ws.Range("M126").Value = -2906.3999
$ws.Range("N126").Value = -46895198

$ws.Range("H134").Value = 3401.9583
$ws.Range("I134").Value = 4083.4707
$ws.Range("J134").Value = 1746.8572
$ws.Range("K134").Value = 12250.4121
$ws.Range("L134").Value = 5240.571599999999
$ws.Range("M134").Value = -9715.4121
$ws.Range("N134").Value = -10310.5716

$ws.Range("H136").Value = 1661.9231
$ws.Range("I136").Value = 925.41174
$ws.Range("J136").Value = 2231.0454
$ws.Range("K136").Value = 2776.23522
$ws.Range("L136").Value = 6693.1362
$ws.Range("M136").Value = -226.23522
$ws.Range("N136").Value = -11793.1362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 995.6667
$ws.Range("I51").Value = 493.5
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 1480.5
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = -1020.5
$ws.Range("N51").Value = -6920

$ws.Range("H68").Value = 1422.7
$ws.Range("I68").Value = 378.375
$ws.Range("K68").Value = 1135.125
$ws.Range("M68").Value = -324.125

$ws.Range("H71").Value = 1422.7
$ws.Range("I71").Value = 378.375
$ws.Range("K71").Value = 3405.375
$ws.Range("M71").Value = 650.625

$ws.Range("H107").Value = 286500
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 286500
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 859500
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -863340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6162.2705
$ws.Range("I70").Value = 6265.269
$ws.Range("J70").Value = 5918.8184
$ws.Range("K70").Value = 6265.269
$ws.Range("L70").Value = 5918.8184
$ws.Range("M70").Value = -5995.269
$ws.Range("N70").Value = -6458.8184

$ws.Range("H73").Value = 6162.2705
$ws.Range("I73").Value = 6265.269
$ws.Range("J73").Value = 5918.8184
$ws.Range("K73").Value = 6265.269
$ws.Range("L73").Value = 5918.8184
$ws.Range("M73").Value = -5329.269
$ws.Range("N73").Value = -7790.8184

$ws.Range("H132").Value = 3285.861
$ws.Range("I132").Value = 3858.111
$ws.Range("J132").Value = 2713.611
$ws.Range("K132").Value = 11574.333
$ws.Range("L132").Value = 8140.833
$ws.Range("M132").Value = -9044.332999999999
$ws.Range("N132").Value = -13200.833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 917.3333
$ws.Range("I31").Value = 784.5714
$ws.Range("J31").Value = 1103.2
$ws.Range("K31").Value = 784.5714
$ws.Range("L31").Value = 1103.2
$ws.Range("M31").Value = -536.5714
$ws.Range("N31").Value = -1599.2

$ws.Range("H40").Value = 1000000000
$ws.Range("I40").Value = 1000000000
$ws.Range("K40").Value = 1000000000
$ws.Range("M40").Value = -999999864

$ws.Range("H93").Value = 62525564
$ws.Range("I93").Value = 50498.25
$ws.Range("J93").Value = 125000620
$ws.Range("K93").Value = 50498.25
$ws.Range("L93").Value = 125000620
$ws.Range("M93").Value = -49250.25
$ws.Range("N93").Value = -125003116

$ws.Range("H122").Value = 11907262
$ws.Range("I122").Value = 35715784
$ws.Range("J122").Value = 3001.25
$ws.Range("K122").Value = 107147352
$ws.Range("L122").Value = 9003.75
$ws.Range("M122").Value = -107144902
$ws.Range("N122").Value = -13903.75

$ws.Range("H132").Value = 57307908
$ws.Range("I132").Value = 76408540
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 229225620
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -229223090
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2220.1072
$ws.Range("I122").Value = 1622.5625
$ws.Range("J122").Value = 3016.8333
$ws.Range("K122").Value = 4867.6875
$ws.Range("L122").Value = 9050.499899999999
$ws.Range("M122").Value = -2417.6875
$ws.Range("N122").Value = -13950.4999

$ws.Range("H136").Value = 4282.1665
$ws.Range("I136").Value = 7899.143
$ws.Range("J136").Value = 1980.4546
$ws.Range("K136").Value = 23697.429
$ws.Range("L136").Value = 5941.3638
$ws.Range("M136").Value = -21147.429
$ws.Range("N136").Value = -11041.3638
